$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) previously listed periods in descending
# order (1801 down to 1607) with the "Valor Mora" column (F) carrying the
# matching amount for each period (1801 => 6571, every other period =>
# 24640). This update re-sorts the statement rows (16-34) into ascending
# period order (1607 up to 1801), keeping each amount attached to its
# period.

$periods = @("1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801")

$startRow = 16
for ($idx = 0; $idx -lt $periods.Length; $idx++) {
    $row = $startRow + $idx
    $period = $periods[$idx]

    if ($period -eq "1801") {
        $amount = 6571
    } else {
        $amount = 24640
    }

    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $amount
}
